$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-132) holds a "Förändrad" (changed) date that was bumped
# by one day, from 2023-09-08 (serial 45177) to 2023-09-09 (serial 45178).
for ($r = 2; $r -le 132; $r++) {
    $ws.Cells.Item($r, 3).Value = 45178
}
